$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1 ("battery"): new row 10 - analog switch ---
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "analog switch"
$ws.Range("C10").Value = "SN74LVC1G3157DCKR"
$ws.Range("D10").Value = 0.0513
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = "https://jlcpcb.com/partdetail/TexasInstruments-SN74LVC1G3157DCKR/C38663"

# --- Section 2 ("Motors"): new row 19 - On/Off switch ---
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = "On/Off switch"
$ws.Range("C19").Value = "SK22D15L5"
$ws.Range("D19").Value = 0.31
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = "https://jlcpcb.com/partdetail/XkbConnection-SK22D15L5/C2884795"

# --- Turn the existing TP4056 URL text in F8 into a real hyperlink, matching
#     the formatting already used by the other link cells (F5, F6) ---
$ws.Hyperlinks.Add($ws.Range("F8"), "https://jlcpcb.com/partdetail/TPOWER-TP4056/C382139")
$ws.Range("F8").Style = $ws.Range("F5").Style

# --- Column C needs to widen to fit the new, longer part number ---
$ws.Columns.Item(3).ColumnWidth = 18.33

# --- Update the active selection left by the editor ---
$ws.Range("A12").Select()
